$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Sending=FAPs,Target=ECs) and row 3 (rename D3 to FAPs) with new figures,
# then add new rows 4-7 to complete the full FAPs/sCs x FAPs/sCs/ECs combination matrix.

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Tnfsf11"
$ws.Cells.Item(2, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.781426666666666
$ws.Cells.Item(2, 8).Value = 5.344279999999999
$ws.Cells.Item(2, 9).Value = 0.9840117515485773
$ws.Cells.Item(2, 10).Value = 0.9840117515485773
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.07580833333333332
$ws.Cells.Item(2, 14).Value = 0.227425
$ws.Cells.Item(2, 15).Value = 0.01036056733223033
$ws.Cells.Item(2, 16).Value = 0.01036056733223033
$ws.Cells.Item(2, 17).Value = 0.1350469865555555
$ws.Cells.Item(2, 18).Value = 1.215422879
$ws.Cells.Item(2, 19).Value = 0.01019492000762494
$ws.Cells.Item(2, 20).Value = 0.01019492000762494

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Tnfsf11"
$ws.Cells.Item(3, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.781426666666666
$ws.Cells.Item(3, 8).Value = 5.344279999999999
$ws.Cells.Item(3, 9).Value = 0.9840117515485773
$ws.Cells.Item(3, 10).Value = 0.9840117515485773
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 6.93109
$ws.Cells.Item(3, 14).Value = 20.79327
$ws.Cells.Item(3, 15).Value = 0.9472576624920083
$ws.Cells.Item(3, 16).Value = 0.9472576624920083
$ws.Cells.Item(3, 17).Value = 12.34722855506667
$ws.Cells.Item(3, 18).Value = 111.1250569956
$ws.Cells.Item(3, 19).Value = 0.9321126716365722
$ws.Cells.Item(3, 20).Value = 0.9321126716365722

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Tnfsf11"
$ws.Cells.Item(4, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.781426666666666
$ws.Cells.Item(4, 8).Value = 5.344279999999999
$ws.Cells.Item(4, 9).Value = 0.9840117515485773
$ws.Cells.Item(4, 10).Value = 0.9840117515485773
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.3101076666666667
$ws.Cells.Item(4, 14).Value = 0.930323
$ws.Cells.Item(4, 15).Value = 0.04238177017576132
$ws.Cells.Item(4, 16).Value = 0.04238177017576132
$ws.Cells.Item(4, 17).Value = 0.5524340669377777
$ws.Cells.Item(4, 18).Value = 4.97190660244
$ws.Cells.Item(4, 19).Value = 0.04170415990438015
$ws.Cells.Item(4, 20).Value = 0.04170415990438015

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Tnfsf11"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.02894466666666666
$ws.Cells.Item(5, 8).Value = 0.08683399999999999
$ws.Cells.Item(5, 9).Value = 0.01598824845142267
$ws.Cells.Item(5, 10).Value = 0.01598824845142267
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.07580833333333332
$ws.Cells.Item(5, 14).Value = 0.227425
$ws.Cells.Item(5, 15).Value = 0.01036056733223033
$ws.Cells.Item(5, 16).Value = 0.01036056733223033
$ws.Cells.Item(5, 17).Value = 0.002194246938888888
$ws.Cells.Item(5, 18).Value = 0.01974822245
$ws.Cells.Item(5, 19).Value = 0.0001656473246053919
$ws.Cells.Item(5, 20).Value = 0.0001656473246053919

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Tnfsf11"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.02894466666666666
$ws.Cells.Item(6, 8).Value = 0.08683399999999999
$ws.Cells.Item(6, 9).Value = 0.01598824845142267
$ws.Cells.Item(6, 10).Value = 0.01598824845142267
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 6.93109
$ws.Cells.Item(6, 14).Value = 20.79327
$ws.Cells.Item(6, 15).Value = 0.9472576624920083
$ws.Cells.Item(6, 16).Value = 0.9472576624920083
$ws.Cells.Item(6, 17).Value = 0.2006180896866666
$ws.Cells.Item(6, 18).Value = 1.80556280718
$ws.Cells.Item(6, 19).Value = 0.01514499085543611
$ws.Cells.Item(6, 20).Value = 0.01514499085543611

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Tnfsf11"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.02894466666666666
$ws.Cells.Item(7, 8).Value = 0.08683399999999999
$ws.Cells.Item(7, 9).Value = 0.01598824845142267
$ws.Cells.Item(7, 10).Value = 0.01598824845142267
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.3101076666666667
$ws.Cells.Item(7, 14).Value = 0.930323
$ws.Cells.Item(7, 15).Value = 0.04238177017576132
$ws.Cells.Item(7, 16).Value = 0.04238177017576132
$ws.Cells.Item(7, 17).Value = 0.008975963042444444
$ws.Cells.Item(7, 18).Value = 0.08078366738199999
$ws.Cells.Item(7, 19).Value = 0.0006776102713811676
$ws.Cells.Item(7, 20).Value = 0.0006776102713811676
